# Update 'F' column values (想去人数 / interest counts) across sheets
# per commit: Update gh-pages to output generated at 456a3b4

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7798
$ws.Range("F5").Value = 7798
$ws.Range("F7").Value = 119
$ws.Range("F8").Value = 2115
$ws.Range("F9").Value = 8590
$ws.Range("F10").Value = 8590
$ws.Range("F13").Value = 89
$ws.Range("F14").Value = 5738
$ws.Range("F16").Value = 2701
$ws.Range("F17").Value = 1191
$ws.Range("F20").Value = 99
$ws.Range("F22").Value = 590
$ws.Range("F23").Value = 23
$ws.Range("F24").Value = 3780
$ws.Range("F29").Value = 94
$ws.Range("F30").Value = 10
$ws.Range("F31").Value = 4104
$ws.Range("F32").Value = 3
$ws.Range("F33").Value = 61
$ws.Range("F35").Value = 380
$ws.Range("F36").Value = 152
$ws.Range("F37").Value = 371
$ws.Range("F38").Value = 1591
$ws.Range("F39").Value = 692
$ws.Range("F42").Value = 3470
$ws.Range("F43").Value = 64
$ws.Range("F45").Value = 30
$ws.Range("F46").Value = 3364
$ws.Range("F48").Value = 2322

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 12
$ws.Range("F10").Value = 6

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7798
$ws.Range("F6").Value = 7798
$ws.Range("F7").Value = 119
$ws.Range("F8").Value = 2115
$ws.Range("F9").Value = 8590
$ws.Range("F10").Value = 8590
$ws.Range("F12").Value = 89
$ws.Range("F13").Value = 5738
$ws.Range("F15").Value = 2701
$ws.Range("F16").Value = 1191
$ws.Range("F19").Value = 99
$ws.Range("F23").Value = 590
$ws.Range("F25").Value = 3780
$ws.Range("F30").Value = 10
$ws.Range("F31").Value = 4104
$ws.Range("F32").Value = 3
$ws.Range("F33").Value = 61
$ws.Range("F34").Value = 380
$ws.Range("F35").Value = 152
$ws.Range("F36").Value = 371
$ws.Range("F38").Value = 1591
$ws.Range("F39").Value = 692
$ws.Range("F43").Value = 3470
$ws.Range("F44").Value = 64
$ws.Range("F46").Value = 30
$ws.Range("F47").Value = 3364
$ws.Range("F48").Value = 2322
$ws.Range("F51").Value = 6

